$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value2 = 111901549
$ws.Range("Q52").Value2 = 477463.5153726833
$ws.Range("R52").Value2 = 7033364.497689161
$ws.Range("AC52").Value2 = "ringhack färska"
$ws.Range("A53").Value2 = 111901547
$ws.Range("B53").Value2 = 56398
$ws.Range("E53").Value2 = 100109
$ws.Range("F53").Value2 = "Tretåig hackspett"
$ws.Range("G53").Value2 = "Picoides tridactylus"
$ws.Range("H53").Value2 = "(Linnaeus, 1758)"
$ws.Range("I53").Value2 = $null
$ws.Range("N53").Value2 = $null
$ws.Range("Q53").Value2 = 477523.7983399219
$ws.Range("R53").Value2 = 7033329.5960882
$ws.Range("AC53").Value2 = "ringhack"
$ws.Range("A54").Value2 = 111901585
$ws.Range("Q54").Value2 = 478338.6535977835
$ws.Range("R54").Value2 = 7035076.017275342
$ws.Range("AC54").Value2 = "ringhack äldre"
$ws.Range("A55").Value2 = 111901584
$ws.Range("Q55").Value2 = 478210.8823140781
$ws.Range("R55").Value2 = 7035067.488268948
$ws.Range("AC55").Value2 = "ringhack"
$ws.Range("A56").Value2 = 111901545
$ws.Range("Q56").Value2 = 477666.6718496145
$ws.Range("R56").Value2 = 7033500.081917049
$ws.Range("AC56").Value2 = "ringhack äldre"
$ws.Range("A57").Value2 = 111901518
$ws.Range("Q57").Value2 = 477673.8480424859
$ws.Range("R57").Value2 = 7033500.479530043
$ws.Range("A58").Value2 = 111901519
$ws.Range("B58").Value2 = 86223
$ws.Range("E58").Value2 = 4412
$ws.Range("F58").Value2 = "Äggvaxskivling"
$ws.Range("G58").Value2 = "Hygrophorus karstenii"
$ws.Range("H58").Value2 = "Sacc. & Cub."
$ws.Range("Q58").Value2 = 477765.0897337386
$ws.Range("R58").Value2 = 7033404.474773662
$ws.Range("AC58").Value2 = $null
$ws.Range("A59").Value2 = 111901544
$ws.Range("Q59").Value2 = 477638.5281090657
$ws.Range("R59").Value2 = 7033514.606123095
$ws.Range("A60").Value2 = 111901551
$ws.Range("B60").Value2 = 56398
$ws.Range("E60").Value2 = 100109
$ws.Range("F60").Value2 = "Tretåig hackspett"
$ws.Range("G60").Value2 = "Picoides tridactylus"
$ws.Range("H60").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q60").Value2 = 477432.586959724
$ws.Range("R60").Value2 = 7033429.191801991
$ws.Range("AC60").Value2 = "ringhack"
$ws.Range("A61").Value2 = 111901546
$ws.Range("Q61").Value2 = 477668.4830064432
$ws.Range("R61").Value2 = 7033374.253324097
$ws.Range("A62").Value2 = 111901587
$ws.Range("B62").Value2 = 56543
$ws.Range("E62").Value2 = 103021
$ws.Range("F62").Value2 = "Talltita"
$ws.Range("G62").Value2 = "Poecile montanus"
$ws.Range("H62").Value2 = "(Conrad von Baldenstein, 1827)"
$ws.Range("I62").Value2 = "'2"
$ws.Range("N62").Value2 = "observerad"
$ws.Range("Q62").Value2 = 477611.096214832
$ws.Range("R62").Value2 = 7033310.625431053
$ws.Range("A63").Value2 = 111901618
$ws.Range("B63").Value2 = 85062
$ws.Range("E63").Value2 = 249278
$ws.Range("F63").Value2 = "Barrviolspindling"
$ws.Range("G63").Value2 = "Cortinarius harcynicus"
$ws.Range("H63").Value2 = "(Pers.) M.M.Moser"
$ws.Range("Q63").Value2 = 477471.0231659907
$ws.Range("R63").Value2 = 7033411.906366638
$ws.Range("AC63").Value2 = $null
$ws.Range("A65").Value2 = 111901619
$ws.Range("B65").Value2 = 85062
$ws.Range("E65").Value2 = 249278
$ws.Range("F65").Value2 = "Barrviolspindling"
$ws.Range("G65").Value2 = "Cortinarius harcynicus"
$ws.Range("H65").Value2 = "(Pers.) M.M.Moser"
$ws.Range("Q65").Value2 = 478523.10621621
$ws.Range("R65").Value2 = 7034650.501669589
$ws.Range("AC65").Value2 = $null
$ws.Range("A66").Value2 = 111901550
$ws.Range("B66").Value2 = 56398
$ws.Range("E66").Value2 = 100109
$ws.Range("F66").Value2 = "Tretåig hackspett"
$ws.Range("G66").Value2 = "Picoides tridactylus"
$ws.Range("H66").Value2 = "(Linnaeus, 1758)"
$ws.Range("Q66").Value2 = 477473.2080285564
$ws.Range("R66").Value2 = 7033403.83150613
$ws.Range("AC66").Value2 = "ringhack äldre"
